$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before the current row 38 (shifts existing row 38 "KFC" and below down by one)
$ws.Rows.Item(38).Insert()

# Populate the new row 38 with the new keyword entry
$ws.Cells.Item(38, 1).Value = "Ingha"
$ws.Cells.Item(38, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
$ws.Cells.Item(38, 3).Value = "Food & Drink"

# Update selection to match the target state
$ws.Range("A39").Select()
